$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UBFO")

# Insert a new column before column D, shifting existing D:K data to E:L.
$ws.Columns("D:D").Insert()

# Copy number formats/styles from the (now shifted) column E into the new column D
# so the new column matches the per-row style (date style for header rows,
# number style for data rows) instead of the generic default.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0


$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 36600
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = "NA"
$ws.Range("D17").Value = 900
$ws.Range("D18").Value = 35700
$ws.Range("D20").Value = -16300
$ws.Range("D21").Value = 20700
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 19300
$ws.Range("D24").Value = 5300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 14000
$ws.Range("D27").Value = 14000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 16300
$ws.Range("D33").Value = 14000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 14000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 220300
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 1900
$ws.Range("D48").Value = 9800
$ws.Range("D49").Value = 4500
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 3200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 933100
$ws.Range("D57").Value = 8000
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 10200
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 823800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 49900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 109200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 14000
$ws.Range("D83").Value = 1400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 11100
$ws.Range("D91").Value = -1000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -11900
$ws.Range("D96").Value = -4700
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 113200
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 112400

# PasteSpecial can materialize stray empty cells on rows that had no cells at
# all in the source range (pure blank separator rows). Clear those back out
# so such rows remain entirely absent, matching the original sheet's layout.
$ws.Range("D36").ClearContents()
$ws.Range("D36").ClearFormats()
$ws.Range("D78").ClearContents()
$ws.Range("D78").ClearFormats()
